$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-7, columns B:E and G (F is unchanged)
$data = @{
    2 = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 0.1496068669990043; E = 0.5333859586016987; G = 3.755628166162433 }
    3 = @{ B = 1.445647641019636;  C = 0.3048912486333797;  D = 0.1496068669990043; E = 0.5333859586016987; G = 2.433531715253719 }
    4 = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 0.1496068669990043; E = 0.5333859586016987; G = 3.755628166162433 }
    5 = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 0.1496068669990043; E = 0.5333859586016987; G = 3.755628166162433 }
    6 = @{ B = 3.272327238179451;  C = 0.04103571897497393; D = 3.223369029078222;  E = 0.5333859586016987; G = 7.070117944834345 }
    7 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
